$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr1 = New-Object "object[,]" 24,4
$arr1[0,0] = 19.3297503425517
$arr1[0,1] = 5.569451790786747
$arr1[0,2] = 13.88280405473222
$arr1[0,3] = 14.11337271956395
$arr1[1,0] = 19.31283106194822
$arr1[1,1] = 5.507570739449757
$arr1[1,2] = 13.88175505347496
$arr1[1,3] = 14.12760739772882
$arr1[2,0] = 19.30783572805636
$arr1[2,1] = 5.467958352679278
$arr1[2,2] = 13.88324079974102
$arr1[2,3] = 14.13755067989542
$arr1[3,0] = 19.30715959359921
$arr1[3,1] = 5.451410753095096
$arr1[3,2] = 13.88438263455601
$arr1[3,3] = 14.14190557527136
$arr1[4,0] = 19.30712950878075
$arr1[4,1] = 5.448638605008783
$arr1[4,2] = 13.88460465149516
$arr1[4,3] = 14.14264700952429
$arr1[5,0] = 19.30782110132124
$arr1[5,1] = 5.467736824771887
$arr1[5,2] = 13.88325402619911
$arr1[5,3] = 14.13760818448812
$arr1[6,0] = 19.32279955936482
$arr1[6,1] = 5.548451666216502
$arr1[6,2] = 13.88200085460966
$arr1[6,3] = 14.11803136351854
$arr1[7,0] = 19.39476900726898
$arr1[7,1] = 5.69384673493118
$arr1[7,2] = 13.89640062435526
$arr1[7,3] = 14.08917061912111
$arr1[8,0] = 19.47328000831292
$arr1[8,1] = 5.792749791212414
$arr1[8,2] = 13.91719112618346
$arr1[8,3] = 14.07375336791664
$arr1[9,0] = 19.51446867436134
$arr1[9,1] = 5.836012603126415
$arr1[9,2] = 13.92884839446381
$arr1[9,3] = 14.06799113724774
$arr1[10,0] = 19.53084308050881
$arr1[10,1] = 5.852145557952025
$arr1[10,2] = 13.93357705736814
$arr1[10,3] = 14.06598858265113
$arr1[11,0] = 19.52728216503339
$arr1[11,1] = 5.848682161775645
$arr1[11,2] = 13.93254471016077
$arr1[11,3] = 14.06641189348842
$arr1[12,0] = 19.51580029215814
$arr1[12,1] = 5.837344882316937
$arr1[12,2] = 13.92923113733507
$arr1[12,3] = 14.06782279145892
$arr1[13,0] = 19.50886820786727
$arr1[13,1] = 5.830367916082546
$arr1[13,2] = 13.92724234793199
$arr1[13,3] = 14.06871036712026
$arr1[14,0] = 19.47069754992644
$arr1[14,1] = 5.789887578743518
$arr1[14,2] = 13.91647339224938
$arr1[14,3] = 14.07415508270854
$arr1[15,0] = 19.44867628306267
$arr1[15,1] = 5.764610598077424
$arr1[15,2] = 13.9104290028405
$arr1[15,3] = 14.07781538475068
$arr1[16,0] = 19.43652610156467
$arr1[16,1] = 5.74990929836829
$arr1[16,2] = 13.90715959589883
$arr1[16,3] = 14.08003848560025
$arr1[17,0] = 19.43250114493713
$arr1[17,1] = 5.744903792334837
$arr1[17,2] = 13.90608826883687
$arr1[17,3] = 14.08081143323564
$arr1[18,0] = 19.45096716274626
$arr1[18,1] = 5.767318217939289
$arr1[18,2] = 13.9110510112415
$arr1[18,3] = 14.07741355145077
$arr1[19,0] = 19.51915178951385
$arr1[19,1] = 5.840681702516527
$arr1[19,2] = 13.93019590064725
$arr1[19,3] = 14.06740350914706
$arr1[20,0] = 19.56823885022645
$arr1[20,1] = 5.887172875620553
$arr1[20,2] = 13.94453915047782
$arr1[20,3] = 14.06190728143179
$arr1[21,0] = 19.54162973816722
$arr1[21,1] = 5.862493251760389
$arr1[21,2] = 13.93671706671872
$arr1[21,3] = 14.06474516981869
$arr1[22,0] = 19.44992986628554
$arr1[22,1] = 5.766094630042653
$arr1[22,2] = 13.91076916071427
$arr1[22,3] = 14.07759485038659
$arr1[23,0] = 19.37076896465189
$arr1[23,1] = 5.655899917147366
$arr1[23,2] = 13.8907063674551
$arr1[23,3] = 14.09596030125861
$ws.Range("B2:E25").Value = $arr1

$arr2 = New-Object "object[,]" 24,2
$arr2[0,0] = 81.47050279165461
$arr2[0,1] = 27.25689947368624
$arr2[1,0] = 80.93854641563044
$arr2[1,1] = 27.19037776352994
$arr2[2,0] = 80.61811373695973
$arr2[2,1] = 27.15171724372863
$arr2[3,0] = 80.48916178843548
$arr2[3,1] = 27.13651673514709
$arr2[4,0] = 80.46784967469203
$arr2[4,1] = 27.13402629184223
$arr2[5,0] = 80.61636796791205
$arr2[5,1] = 27.15150999545724
$arr2[6,0] = 81.2858423616495
$arr2[6,1] = 27.23351244391575
$arr2[7,0] = 82.64459933908637
$arr2[7,1] = 27.41142822816109
$arr2[8,0] = 83.66644095565979
$arr2[8,1] = 27.55226181828962
$arr2[9,0] = 84.13543338756133
$arr2[9,1] = 27.61845400739077
$arr2[10,0] = 84.313539380809
$arr2[10,1] = 27.64381797705997
$arr2[11,0] = 84.2751598594465
$arr2[11,1] = 27.63834224042964
$arr2[12,0] = 84.15007658853747
$arr2[12,1] = 27.62053477885885
$arr2[13,0] = 84.07352308699694
$arr2[13,1] = 27.60966583431941
$arr2[14,0] = 83.63586722106537
$arr2[14,1] = 27.54797817077764
$arr2[15,0] = 83.36838044571032
$arr2[15,1] = 27.51067427790116
$arr2[16,0] = 83.21492526119201
$arr2[16,1] = 27.48941866348632
$arr2[17,0] = 83.1630386773797
$arr2[17,1] = 27.48225654799769
$arr2[18,0] = 83.39681450687443
$arr2[18,1] = 27.51462463194614
$arr2[19,0] = 84.18680347927737
$arr2[19,1] = 27.62575722305062
$arr2[20,0] = 84.70603058039245
$arr2[20,1] = 27.70012458846127
$arr2[21,0] = 84.42867122695657
$arr2[21,1] = 27.66027701557451
$arr2[22,0] = 83.3839584439831
$arr2[22,1] = 27.51283808167387
$arr2[23,0] = 82.2725518482705
$arr2[23,1] = 27.36149620462859
$ws.Range("G2:H25").Value = $arr2

$arr3 = New-Object "object[,]" 24,4
$arr3[0,0] = 8.547109886358177
$arr3[0,1] = 14.98036617211361
$arr3[0,2] = 12.56499159740429
$arr3[0,3] = 19.0298370435161
$arr3[1,0] = 8.548235606665781
$arr3[1,1] = 14.98763983043296
$arr3[1,2] = 12.58777098321887
$arr3[1,3] = 19.06426905374
$arr3[2,0] = 8.548970536547303
$arr3[2,1] = 14.99635513786956
$arr3[2,2] = 12.60334933111026
$arr3[2,3] = 19.08875455146821
$arr3[3,0] = 8.549281039748521
$arr3[3,1] = 15.0009753467167
$arr3[3,2] = 12.61009820834672
$arr3[3,3] = 19.09957351615929
$arr3[4,0] = 8.549333264228142
$arr3[4,1] = 15.00180706444577
$arr3[4,2] = 12.61124305710528
$arr3[4,3] = 19.10142078517563
$arr3[5,0] = 8.548974679486619
$arr3[5,1] = 14.99641312100877
$arr3[5,2] = 12.60343872650772
$arr3[5,3] = 19.08889705500956
$arr3[6,0] = 8.547488967711569
$arr3[6,1] = 14.98199246798045
$arr3[6,2] = 12.57251576954199
$arr3[6,3] = 19.04101517849546
$arr3[7,0] = 8.544921755483815
$arr3[7,1] = 14.98739684551558
$arr3[7,2] = 12.52449103121321
$arr3[7,3] = 18.97364885786973
$arr3[8,0] = 8.543245690375882
$arr3[8,1] = 15.01181843861911
$arr3[8,2] = 12.49687804601807
$arr3[8,3] = 18.94031794642191
$arr3[9,0] = 8.542528605203181
$arr3[9,1] = 15.02733744707352
$arr3[9,2] = 12.4859772716574
$arr3[9,3] = 18.92865975742615
$arr3[10,0] = 8.542263571916649
$arr3[10,1] = 15.03384433609729
$arr3[10,2] = 12.48208779336517
$arr3[10,3] = 18.92474833159933
$arr3[11,0] = 8.542320362175362
$arr3[11,1] = 15.03241500336569
$arr3[11,2] = 12.4829148650095
$arr3[11,3] = 18.92556835529199
$arr3[12,0] = 8.542506670326919
$arr3[12,1] = 15.02786016789038
$arr3[12,2] = 12.48565250609428
$arr3[12,3] = 18.9283278797568
$arr3[13,0] = 8.542621637059005
$arr3[13,1] = 15.02515213959733
$arr3[13,2] = 12.48736042666842
$arr3[13,3] = 18.93008368862352
$arr3[14,0] = 8.543293465396074
$arr3[14,1] = 15.01089264759296
$arr3[14,2] = 12.49762381984253
$arr3[14,3] = 18.94115030029476
$arr3[15,0] = 8.543717220689931
$arr3[15,1] = 15.0032717344433
$arr3[15,2] = 12.50434511978641
$arr3[15,3] = 18.94883642451906
$arr3[16,0] = 8.543965224225991
$arr3[16,1] = 14.99930370823296
$arr3[16,2] = 12.50836735272239
$arr3[16,3] = 18.9535871596354
$arr3[17,0] = 8.544049927834864
$arr3[17,1] = 14.99803164208557
$arr3[17,2] = 12.50975607049678
$arr3[17,3] = 18.95525235056963
$arr3[18,0] = 8.543671669302114
$arr3[18,1] = 15.00404003561774
$arr3[18,2] = 12.50361345014764
$arr3[18,3] = 18.94798408627712
$arr3[19,0] = 8.542451770487267
$arr3[19,1] = 15.02918096520592
$arr3[19,2] = 12.48484192729477
$arr3[19,3] = 18.92750368745306
$arr3[20,0] = 8.541692443562159
$arr3[20,1] = 15.04928287182337
$arr3[20,2] = 12.4739630882318
$arr3[20,3] = 18.91705165087063
$arr3[21,0] = 8.542094242364939
$arr3[21,1] = 15.03821969138265
$arr3[21,2] = 12.47964232506611
$arr3[21,3] = 18.92236196733676
$arr3[22,0] = 8.543692249438694
$arr3[22,1] = 15.0036913988701
$arr3[22,2] = 12.50394374560405
$arr3[22,3] = 18.94836839458634
$arr3[23,0] = 8.545579291730409
$arr3[23,1] = 14.98233469394979
$arr3[23,2] = 12.5361344115246
$arr3[23,3] = 18.98903353682016
$ws.Range("J2:M25").Value = $arr3

Write-Host "applied loading_percent updates"